$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 150
$ws.Range("E3").Value = 1150
$ws.Range("E4").Value = 1515
$ws.Range("E5").Value = 10150
$ws.Range("E6").Value = 150
$ws.Range("E7").Value = 1150
$ws.Range("E8").Value = 150
$ws.Range("E9").Value = 115

$ws.Range("A10").Value = "Super Vaca Lechera"
$ws.Range("B10").Value = "Planta del Prado"
$ws.Range("C10").Value = 160
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 150
